$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as plain text (some use "." as a thousands
# separator, e.g. "42.488.50"). For values that otherwise look like a plain
# number, force text entry with a leading apostrophe (like typing into
# Excel) so the cell keeps its string type instead of becoming numeric, then
# reset the style so no left-over "quote prefix" formatting lingers on the
# cell.
$ws.Range("D2").Value = "42.488.50"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.183.41"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  -0.14%  "
$cell = $ws.Range("D5")
$cell.Value = "'250.87"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$cell = $ws.Range("D6")
$cell.Value = "'0.614"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "
$cell = $ws.Range("D7")
$cell.Value = "'75.09"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -0.06%  "
$cell = $ws.Range("D9")
$cell.Value = "'0.583"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -5.35%  "
$cell = $ws.Range("D10")
$cell.Value = "'40.25"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "
$cell = $ws.Range("D11")
$cell.Value = "'0.0910"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.21%  "
$cell = $ws.Range("D12")
$cell.Value = "'0.101"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "
$cell = $ws.Range("D13")
$cell.Value = "'6.78"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "2.514.18"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "2.186.05"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("E17").Value = "  -5.46%  "
$ws.Range("D18").Value = "42.390.76"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -3.18%  "
$cell = $ws.Range("D20")
$cell.Value = "'70.72"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$cell = $ws.Range("D21")
$cell.Value = "'5.86"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.19%  "
$cell = $ws.Range("D22")
$cell.Value = "'225.54"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  -10.29%  "
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("E25").Value = "  -0.02%  "
$cell = $ws.Range("D26")
$cell.Value = "'10.43"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -5.00%  "
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +0.18%  "
$cell = $ws.Range("D29")
$cell.Value = "'37.64"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("E30").Value = "  -4.56%  "
$cell = $ws.Range("D31")
$cell.Value = "'172.38"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -1.67%  "
$cell = $ws.Range("D33")
$cell.Value = "'0.0814"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  -3.48%  "
$cell = $ws.Range("D38")
$cell.Value = "'0.0334"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "
$cell = $ws.Range("D39")
$cell.Value = "'11.95"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -8.63%  "
$ws.Range("E40").Value = "  -3.37%  "
$cell = $ws.Range("D41")
$cell.Value = "'2.57"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +10.94%  "
$ws.Range("E42").Value = "  -7.78%  "
$ws.Range("E43").Value = "  -2.77%  "
$cell = $ws.Range("D44")
$cell.Value = "'58.71"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.33%  "
$cell = $ws.Range("D45")
$cell.Value = "'101.37"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.79%  "
$cell = $ws.Range("D47")
$cell.Value = "'0.459"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -0.94%  "
